$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update B2, E2, F2 values (prodctname -> Baclpanle, 55555 -> 600, 999 -> 300)
$ws.Range("B2").Value = "Baclpanle"
$ws.Range("E2").Value = "600"
$ws.Range("F2").Value = "300"

# Clear D2, G2, H2, J2, K2, P2 entirely (cells no longer present)
$ws.Range("D2").Clear()
$ws.Range("G2").Clear()
$ws.Range("H2").Clear()
$ws.Range("J2").Clear()
$ws.Range("K2").Clear()
$ws.Range("P2").Clear()

# Update the view selection to P2 (matches the scrolled/selected state in the target workbook)
$ws.Range("P2").Select()
